# Apply crypto "Price" (D) and "Volume(1h)" (E) updates for rows 2-51 per the commit diff.
# A leading apostrophe forces Excel to store dotted/zero-padded numeric-looking strings
# (e.g. "0.1200", "7.330") as text instead of silently coercing them to numbers and
# losing their exact formatting - matching how these cells were authored originally.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.894.03"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.632.25"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'209.59"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").Value = "'0.5206"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.2564"
$ws.Range("E8").Value = "  -3.45%  "
$ws.Range("D9").Value = "'0.06228"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "'20.23"
$ws.Range("D11").Value = "'0.07556"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "1.664.09"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Value = "'4.351"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").Value = "1.858.79"
$ws.Range("E14").Value = "  -2.50%  "
$ws.Range("D15").Value = "'0.5409"
$ws.Range("E15").Value = "  -4.08%  "
$ws.Range("D16").Value = "0.0₅7911"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "'64.41"
$ws.Range("E17").Value = "  -3.79%  "
$ws.Range("D18").Value = "25.891.23"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'4.607"
$ws.Range("E20").Value = "  -4.53%  "
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").Value = "'9.993"
$ws.Range("E22").Value = "  -4.17%  "
$ws.Range("D23").Value = "'6.053"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'145.61"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").Value = "'0.1200"
$ws.Range("E26").Value = "  -3.88%  "
$ws.Range("D27").Value = "'7.330"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("D28").Value = "'15.47"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("D29").Value = "'1.373"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").Value = "'0.05928"
$ws.Range("E30").Value = "  -4.34%  "
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("D32").Value = "'3.344"
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("D33").Value = "'3.343"
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("D34").Value = "'1.603"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").Value = "'0.9686"
$ws.Range("D36").Value = "'2.382"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Value = "'2.733"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").Value = "'0.5794"
$ws.Range("E38").Value = "  -4.44%  "
$ws.Range("D39").Value = "'0.01590"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").Value = "'0.8386"
$ws.Range("E41").Value = "  -3.48%  "
$ws.Range("D42").Value = "'5.650"
$ws.Range("E42").Value = "  -7.34%  "
$ws.Range("D43").Value = "1.021.03"
$ws.Range("E43").Value = "  -5.81%  "
$ws.Range("D44").Value = "'99.53"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").Value = "1.784.17"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("E46").Value = "  -4.32%  "
$ws.Range("D47").Value = "'0.9998"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "'54.23"
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("D49").Value = "'7.967"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "'0.4227"
